$wb = $excel.ActiveWorkbook

# --- "Simple Fields" sheet: Currency SGD -> USD ---
$ws1 = $wb.Worksheets.Item("Simple Fields")
$ws1.Range("I2").Value = "USD"

# --- "Simple Fields - Formatted" sheet: Currency SGD -> USD ---
$ws2 = $wb.Worksheets.Item("Simple Fields - Formatted")
$ws2.Range("I2").Value = "USD"

# --- "Items" sheet: item description corrections ---
$ws3 = $wb.Worksheets.Item("Items")
$ws3.Range("A2").Value = "green onion Pancakes (1)"
$ws3.Range("A3").Value = "Pan Fried Leek Dumplings IAT (2)"
$ws3.Range("A4").Value = "Pork Xiao Long Bao(10) A¥R]J¿E(10)"
$ws3.Range("A5").Value = "Q-BAO (5) WEEL (5)"
$ws3.Range("A6").Value = "Chicken potstickers KR$55(6)"
$ws3.Range("A7").Value = "Tomato Mushroom Steamed dumpli pEiAINABUNXA (6)"
$ws3.Range("A8").Value = "Zucchini shrimp dumplings A/LC"
$ws3.Range("A9").Value = "beef stew nodle soup (Non Spicy 0H#4PJB(TY)"
$ws3.Range("A10").Value = "dandan noodle INCMM"
$ws3.Range("A11").Value = "banana naan bread BATAI"
$ws3.Range("A12").Value = "house made plum juice"

# --- "Items - Formatted" sheet: item description corrections ---
$ws4 = $wb.Worksheets.Item("Items - Formatted")
$ws4.Range("A2").Value = "green onion Pancakes (1)"
$ws4.Range("A3").Value = "Pan Fried Leek Dumplings IAT (2)"
$ws4.Range("A4").Value = "Pork Xiao Long Bao(10) A¥R]J¿E(10)"
$ws4.Range("A5").Value = "Q-BAO (5) WEEL (5)"
$ws4.Range("A6").Value = "Chicken potstickers KR$55(6)"
$ws4.Range("A7").Value = "Tomato Mushroom Steamed dumpli pEiAINABUNXA (6)"
$ws4.Range("A8").Value = "Zucchini shrimp dumplings A/LC"
$ws4.Range("A9").Value = "beef stew nodle soup (Non Spicy 0H#4PJB(TY)"
$ws4.Range("A10").Value = "dandan noodle INCMM"
$ws4.Range("A11").Value = "banana naan bread BATAI"
$ws4.Range("A12").Value = "house made plum juice"
